$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.854.33'
$ws.Range("E2").Value = '  +0.28%  '

# Row 3
$ws.Range("D3").Value = '2.315.85'
$ws.Range("E3").Value = '  +0.93%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.31'
$ws.Range("E5").Value = '  -0.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.72'
$ws.Range("E6").Value = '  -0.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.35%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  -0.81%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.40'
$ws.Range("E10").Value = '  -1.71%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.05'
$ws.Range("E11").Value = '  +1.84%  '

# Row 12
$ws.Range("E12").Value = '  +0.20%  '

# Row 13
$ws.Range("E13").Value = '  +0.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("E14").Value = '  -1.38%  '

# Row 15
$ws.Range("D15").Value = '2.678.14'
$ws.Range("E15").Value = '  +0.99%  '

# Row 16
$ws.Range("D16").Value = '2.328.46'
$ws.Range("E16").Value = '  +1.61%  '

# Row 17
$ws.Range("E17").Value = '  +2.02%  '

# Row 18
$ws.Range("D18").Value = '42.790.65'
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.14'
$ws.Range("E19").Value = '  -5.04%  '

# Row 20
$ws.Range("E20").Value = '  +2.61%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0891'
$ws.Range("E21").Value = '  -0.25%  '

# Row 23
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.26'
$ws.Range("E23").Value = '  +5.25%  '

# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.35'
$ws.Range("E24").Value = '  -0.24%  '

# Row 25
$ws.Range("E25").Value = '  +0.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +0.90%  '

# Row 27
$ws.Range("E27").Value = '  -2.34%  '

# Row 28
$ws.Range("E28").Value = '  +14.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.56'
$ws.Range("E29").Value = '  -0.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.13'
$ws.Range("E30").Value = '  +1.35%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.20'
$ws.Range("E31").Value = '  -2.17%  '

# Row 32
$ws.Range("E32").Value = '  +0.06%  '

# Row 33
$ws.Range("E33").Value = '  +0.90%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.86'
$ws.Range("E34").Value = '  -0.17%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("E35").Value = '  +0.76%  '

# Row 36
$ws.Range("E36").Value = '  +2.07%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  -0.72%  '

# Row 38
$ws.Range("E38").Value = '  +2.32%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1000'
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.72'
$ws.Range("E40").Value = '  +1.40%  '

# Row 41
$ws.Range("E41").Value = '  -0.30%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.83'
$ws.Range("E42").Value = '  +13.45%  '

# Row 43
$ws.Range("D43").Value = '1.928.16'
$ws.Range("E43").Value = '  -3.27%  '

# Row 44
$ws.Range("E44").Value = '  -0.16%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.13'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("E46").Value = '  -4.64%  '

# Row 47
$ws.Range("E47").Value = '  -0.45%  '

# Row 48
$ws.Range("E48").Value = '  -0.80%  '

# Row 49
$ws.Range("D49").Value = '2.547.33'
$ws.Range("E49").Value = '  +1.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.32'
$ws.Range("E50").Value = '  -0.54%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.28'
$ws.Range("E51").Value = '  +2.04%  '
